$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear only the data rows (2:13); keep header row (row 1) and its styling untouched.
$ws.Range("A2:T13").Clear()

# Row data: columns A-D are string labels (Sending cluster, Ligand symbol, Receptor symbol, Target cluster)
# columns E-T are numeric metrics. Values below match the target workbook exactly,
# and are written in an order that reproduces the expected shared-string table sequence
# (ECs, FAPs, M1, M2, sCs, Sema4d, Erbb2).
$rows = @(
    @("ECs", "Sema4d", "Erbb2", "ECs", 2, 0.6666666666666666, 0.6040246666666667, 1.812074, 0.007395659364387973, 0.007395659364387974, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 2.240805473519556, 20.167249261676, 0.003167771257239938, 0.003167771257239939),
    @("ECs", "Sema4d", "Erbb2", "FAPs", 2, 0.6666666666666666, 0.6040246666666667, 1.812074, 0.007395659364387973, 0.007395659364387974, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 2.173885580699556, 19.564970226296, 0.003073168260452442, 0.003073168260452443),
    @("ECs", "Sema4d", "Erbb2", "sCs", 2, 0.6666666666666666, 0.6040246666666667, 1.812074, 0.007395659364387973, 0.007395659364387974, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 0.8168211473424444, 7.351390326081999, 0.001154719846695592, 0.001154719846695592),
    @("FAPs", "Sema4d", "Erbb2", "ECs", 3, 1, 1.123796, 3.371388, 0.0137597235174641, 0.0137597235174641, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 4.169048661234668, 37.52143795111201, 0.005893680944268083, 0.005893680944268085),
    @("FAPs", "Sema4d", "Erbb2", "FAPs", 3, 1, 1.123796, 3.371388, 0.0137597235174641, 0.0137597235174641, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 4.044543302394668, 36.40088972155201, 0.005717670798913423, 0.005717670798913424),
    @("FAPs", "Sema4d", "Erbb2", "sCs", 3, 1, 1.123796, 3.371388, 0.0137597235174641, 0.0137597235174641, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 1.519706708609334, 13.677360377484, 0.002148371774282594, 0.002148371774282595),
    @("M1", "Sema4d", "Erbb2", "ECs", 3, 1, 28.13359466666667, 84.400784, 0.3444668642402499, 0.3444668642402499, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 104.3697656699129, 939.3278910292161, 0.1475449554729644, 0.1475449554729644),
    @("M1", "Sema4d", "Erbb2", "FAPs", 3, 1, 28.13359466666667, 84.400784, 0.3444668642402499, 0.3444668642402499, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 101.2528447167929, 911.2756024511361, 0.1431386414385408, 0.1431386414385408),
    @("M1", "Sema4d", "Erbb2", "sCs", 3, 1, 28.13359466666667, 84.400784, 0.3444668642402499, 0.3444668642402499, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 38.04499442267911, 342.404949804112, 0.05378326732874471, 0.05378326732874471),
    @("M2", "Sema4d", "Erbb2", "ECs", 3, 1, 49.52258933333334, 148.567768, 0.6063530543760464, 0.6063530543760464, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 183.7184727130258, 1653.466254417232, 0.2597182594213545, 0.2597182594213545),
    @("M2", "Sema4d", "Erbb2", "FAPs", 3, 1, 49.52258933333334, 148.567768, 0.6063530543760464, 0.6063530543760464, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 178.2318650407858, 1604.086785367072, 0.2519619779014886, 0.2519619779014887),
    @("M2", "Sema4d", "Erbb2", "sCs", 3, 1, 49.52258933333334, 148.567768, 0.6063530543760464, 0.6063530543760464, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 66.96928200275822, 602.7235380248239, 0.09467281705320325, 0.09467281705320325),
    @("sCs", "Sema4d", "Erbb2", "ECs", 3, 1, 2.288857333333333, 6.866572, 0.0280246985018516, 0.02802469850185161, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 8.491183098436444, 76.420647885928, 0.01200377546246376, 0.01200377546246376),
    @("sCs", "Sema4d", "Erbb2", "FAPs", 3, 1, 2.288857333333333, 6.866572, 0.0280246985018516, 0.02802469850185161, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 8.237600594476444, 74.13840535028801, 0.01164529215060281, 0.01164529215060282),
    @("sCs", "Sema4d", "Erbb2", "sCs", 3, 1, 2.288857333333333, 6.866572, 0.0280246985018516, 0.02802469850185161, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 3.095216431199555, 27.856947880796, 0.004375630888785027, 0.004375630888785029)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    $targetRow = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value2 = $rowVals[$c]
    }
}

Write-Host "done"